$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (currently holding 005161568 / MARCO),
# shifting that row -- and every row below it -- down by one.
$ws.Rows.Item(5).Insert()

# Column A holds account numbers that must stay as text (leading zeros
# need to survive). Assigning a numeric-looking string straight to
# .Value makes Excel coerce it to a number, and prefixing it with a
# leading apostrophe stores it as text but leaves a "quote prefix" cell
# style behind that the source data never had. Instead, build the text
# value with a formula (formula results are never auto-coerced to a
# number), then Copy/PasteSpecial (values only) it into place -- this
# drops the formula and leaves a plain text cell with no extra style.
$scratch = $ws.Cells.Item(1048576, 16384)
$scratch.Formula = "=""004350197"""
$scratch.Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item(5, 2).Value = "GISELA"
$ws.Cells.Item(5, 3).Value = 8020.41
